# Fruta / hortaliza, semanal
# Insert 3 new weekly price-report rows (Navel Late) above the existing
# row 432 in "Vega Modelo de Temuco" / Naranja data, pushing the old
# rows 432-445 down to 435-448.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 432; everything from 432 down shifts by 3.
$ws.Rows("432:434").Insert()

# Row 432: Navel Late / Primera
$ws.Cells.Item(432, 1).Value = 10
$ws.Cells.Item(432, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(432, 3).Value = "La Araucanía"
$ws.Cells.Item(432, 4).Value = 44448
$ws.Cells.Item(432, 5).Value = 9
$ws.Cells.Item(432, 6).Value = "Fruta"
$ws.Cells.Item(432, 7).Value = 100102
$ws.Cells.Item(432, 8).Value = "Cítricos"
$ws.Cells.Item(432, 9).Value = 100102005
$ws.Cells.Item(432, 10).Value = "Naranja"
$ws.Cells.Item(432, 11).Value = "Navel Late"
$ws.Cells.Item(432, 12).Value = "Primera"
$ws.Cells.Item(432, 13).Value = 400
$ws.Cells.Item(432, 14).Value = 8000
$ws.Cells.Item(432, 15).Value = 8000
$ws.Cells.Item(432, 16).Value = 8000
$ws.Cells.Item(432, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(432, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(432, 19).Value = 533
$ws.Cells.Item(432, 20).Value = 15

# Row 433: Navel Late / Primera (bins)
$ws.Cells.Item(433, 1).Value = 10
$ws.Cells.Item(433, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(433, 3).Value = "La Araucanía"
$ws.Cells.Item(433, 4).Value = 44448
$ws.Cells.Item(433, 5).Value = 9
$ws.Cells.Item(433, 6).Value = "Fruta"
$ws.Cells.Item(433, 7).Value = 100102
$ws.Cells.Item(433, 8).Value = "Cítricos"
$ws.Cells.Item(433, 9).Value = 100102005
$ws.Cells.Item(433, 10).Value = "Naranja"
$ws.Cells.Item(433, 11).Value = "Navel Late"
$ws.Cells.Item(433, 12).Value = "Primera"
$ws.Cells.Item(433, 13).Value = 15
$ws.Cells.Item(433, 14).Value = 140000
$ws.Cells.Item(433, 15).Value = 150000
$ws.Cells.Item(433, 16).Value = 143333
$ws.Cells.Item(433, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(433, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(433, 19).Value = 358
$ws.Cells.Item(433, 20).Value = 400

# Row 434: Navel Late / Segunda
$ws.Cells.Item(434, 1).Value = 10
$ws.Cells.Item(434, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(434, 3).Value = "La Araucanía"
$ws.Cells.Item(434, 4).Value = 44448
$ws.Cells.Item(434, 5).Value = 9
$ws.Cells.Item(434, 6).Value = "Fruta"
$ws.Cells.Item(434, 7).Value = 100102
$ws.Cells.Item(434, 8).Value = "Cítricos"
$ws.Cells.Item(434, 9).Value = 100102005
$ws.Cells.Item(434, 10).Value = "Naranja"
$ws.Cells.Item(434, 11).Value = "Navel Late"
$ws.Cells.Item(434, 12).Value = "Segunda"
$ws.Cells.Item(434, 13).Value = 120
$ws.Cells.Item(434, 14).Value = 3500
$ws.Cells.Item(434, 15).Value = 4000
$ws.Cells.Item(434, 16).Value = 3667
$ws.Cells.Item(434, 17).Value = "$/malla 15 kilos"
$ws.Cells.Item(434, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(434, 19).Value = 244
$ws.Cells.Item(434, 20).Value = 15
